# Daily attendance processing - 2026-01-29 23:41:47
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (G) wherever the old "System, <email>" ordering
# is found, so it reads "<email>, System" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
